$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")
$ws.Activate()

# The "wbGetDataTask_Type" / "wbGetDataTask_SuppressSuccessful" pair (rows 23-24)
# is being removed -- this was the service-run-once GetData task row pair that
# is no longer part of the config list. Deleting the two rows shifts every
# row below them up by two, which is exactly what the target workbook shows.
$ws.Rows.Item(23).Resize(2).EntireRow.Delete() | Out-Null

# Excel leaves the selection wherever the user last clicked before saving;
# the target file shows B31 selected.
$ws.Range("B31").Select() | Out-Null
